$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: section title (merged A32:C32) ---
$ws.Range("A32").Value = "Railsync Driving load sense"
[void]$ws.Range("A32:C32").Merge()

# --- Row 33: R26 ---
$ws.Range("A33").Value = "R26"
$ws.Range("B33").Value = 3000
$ws.Range("C33").Formula = '=B33*0.95'
$ws.Range("D33").Formula = '=B33*1.05'

# --- Row 34: R27 ---
$ws.Range("A34").Value = "R27"
$ws.Range("B34").Value = 22000
$ws.Range("C34").Formula = '=B34*1.05'
$ws.Range("D34").Formula = '=B34*0.95'

# --- Row 35: Ratio ---
$ws.Range("A35").Value = "Ratio"
$ws.Range("B35").Formula = '=B34/B33'
$ws.Range("C35").Formula = '=C34/C33'
$ws.Range("D35").Formula = '=D34/D33'

# --- Row 36: Rtotal ---
$ws.Range("A36").Value = "Rtotal"
$ws.Range("B36").Formula = '=B34+B33'
$ws.Range("C36").Formula = '=C34+C33'
$ws.Range("D36").Formula = '=D34+D33'

# --- Row 37: headers ---
$ws.Range("A37").Value = "RS Volts"
$ws.Range("B37").Value = "Sense Volts"
$ws.Range("C37").Value = "`"-R26 +R27”"
$ws.Range("D37").Value = "“+R26 -R27”"
$ws.Range("E37").Value = "Count"

# --- Row 38 ---
$ws.Range("A38").Value = 4
$ws.Range("B38").Formula = '=$A38/B$36*B$33'
$ws.Range("C38").Formula = '=$A38/C$36*C$33'
$ws.Range("D38").Formula = '=$A38/D$36*D$33'
$ws.Range("E38").Formula = '=ROUND($C$14*B38,0)'

# --- Row 39 ---
$ws.Range("A39").Value = 5
$ws.Range("B39").Formula = '=A39/$B$36*$B$33'
$ws.Range("C39").Formula = '=$A39/C$36*C$33'
$ws.Range("D39").Formula = '=$A39/D$36*D$33'
$ws.Range("E39").Formula = '=ROUND($C$14*B39,0)'

# --- Row 40 ---
$ws.Range("A40").Value = 7
$ws.Range("B40").Formula = '=A40/$B$36*$B$33'
$ws.Range("C40").Formula = '=$A40/C$36*C$33'
$ws.Range("D40").Formula = '=$A40/D$36*D$33'
$ws.Range("E40").Formula = '=ROUND($C$14*B40,0)'

# --- Row 41 ---
$ws.Range("A41").Value = 12
$ws.Range("B41").Formula = '=A41/$B$36*$B$33'
$ws.Range("C41").Formula = '=$A41/C$36*C$33'
$ws.Range("D41").Formula = '=$A41/D$36*D$33'
$ws.Range("E41").Formula = '=ROUND($C$14*B41,0)'

# --- Row 42 ---
$ws.Range("A42").Value = 15
$ws.Range("B42").Formula = '=A42/$B$36*$B$33'
$ws.Range("C42").Formula = '=$A42/C$36*C$33'
$ws.Range("D42").Formula = '=$A42/D$36*D$33'
$ws.Range("E42").Formula = '=ROUND($C$14*B42,0)'

# --- Row 43 ---
$ws.Range("A43").Value = 24
$ws.Range("B43").Formula = '=A43/$B$36*$B$33'
$ws.Range("C43").Formula = '=$A43/C$36*C$33'
$ws.Range("D43").Formula = '=$A43/D$36*D$33'
$ws.Range("E43").Formula = '=ROUND($C$14*B43,0)'

# --- Selection / active cell ---
[void]$ws.Range("E38").Select()

Write-Host "done"
